$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new hospital ("Meander Medisch Centrum") needs to be inserted into the
# alphabetically-sorted list at row 30 (between "Martini Ziekenhuis" and
# "Medisch Centrum Leeuwarden"), pushing every row below it down by one.
$ws.Rows.Item(30).Insert()

$ws.Cells.Item(30, 1).Value = "Meander Medisch Centrum"
$ws.Cells.Item(30, 2).Value = "definitief en/of vastgesteld RvB"
